$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 9.5
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 5
$ws.Range("C14").Value = 5

$ws.Range("D15").Value = "** quand le menu grossi, le header size grossi avec, contact us, Mailing List  et about us on selement besoins d'apparaitre une fois dans la page"
$ws.Range("D5").Value = "*bonne addition avec le boutton add to cart, un icone de panier pourrais être pratique si on veux que le consomateur puissse acheter le produit"
$ws.Range("D6").Value = "* style autour du menu fait bouger les choses avec un dropdown"
$ws.Range("D13").Value = "* quelques commentaires en début de page, mais subsidues vers la fin"

$ws.Range("C15").Select()
